$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.276.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.50%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.845.64"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.24%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.9987"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.08%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'241.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.09%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.6737"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.80%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.9999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.01%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.07442"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.32%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.2954"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.19%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'22.91"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -1.16%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07721"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.69%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.838.30"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.01%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'5.008"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.12%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.6733"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.53%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'86.22"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.73%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'6.139"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.41%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'29.242.47"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.38%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.000008333"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.96%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'228.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.44%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'12.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.06%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'1.000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.07%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'7.200"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -2.80%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.9999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.00%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'160.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.45%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'8.725"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.46%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.1406"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -3.54%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'18.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.39%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'1.511"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.36%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'4.180"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -2.13%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'4.076"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -1.77%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +0.04%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.05312"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +2.49%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.7610"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.69%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.875"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.36%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.138"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.12%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'2.676"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Value = "'1.329.14"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.82%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.01805"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -1.76%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -0.14%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.9199"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.80%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'6.004"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +4.07%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +0.21%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'103.43"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.47%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.08172"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +12.00%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "'RocketPoolETH"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'1.991.48"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.43%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'BabyDogeCoin"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000124"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.95%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.5167"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.54%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'1.781"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.25%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'63.88"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.77%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'9.194"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -3.81%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.05955"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.19%  "
$ws.Range("E51").Style = "Normal"
